$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Marking" row total right-answer value (B11): 3 -> 5
$ws.Range("B11").Value = 5

# Update "Total" row value (B12): 45 -> 75
$ws.Range("B12").Value = 75

# Update correct/total marks summary text (E12): "44/84" -> "75/140"
$ws.Range("E12").Value = "75/140"
